$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "status_label" column (string version of the emoji "statut" column)
# before the existing NCTId column; this shifts NCTId..results from B..I to C..J.
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Cells.Item(1, 2).Value2 = "status_label"

# The source data for rows 7 and 8 was also reordered (NCT04634318/REHABCOVID now
# comes before NCT05237050/SONOMYAL); swap every column that differs between the two rows.
$swapCols = @("A","C","F","G","H","I","J")
foreach ($c in $swapCols) {
    $v7 = $ws.Range("$c`7").Value2
    $v8 = $ws.Range("$c`8").Value2
    $ws.Range("$c`7").Value2 = $v8
    $ws.Range("$c`8").Value2 = $v7
}

# Fill the new status_label column with the French word matching each row's emoji status
$ws.Range("B2").Value2 = "rouge"
$ws.Range("B3").Value2 = "vert"
$ws.Range("B4").Value2 = "rouge"
$ws.Range("B5").Value2 = "rouge"
$ws.Range("B6").Value2 = "rouge"
$ws.Range("B7").Value2 = "orange"
$ws.Range("B8").Value2 = "rouge"
$ws.Range("B9").Value2 = "rouge"
$ws.Range("B10").Value2 = "rouge"
$ws.Range("B11").Value2 = "rouge"
$ws.Range("B12").Value2 = "rouge"
$ws.Range("B13").Value2 = "rouge"
$ws.Range("B14").Value2 = "orange"
$ws.Range("B15").Value2 = "rouge"
$ws.Range("B16").Value2 = "rouge"

